$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = '2011-2014'
$ws.Range("B26").Value = '**島根県保健環境科学研究所** <br> [ウズラ卵が原因食品と推定された _Salmonella enterica_ serovar 4，［5］，12：i：一による食中毒の発生とウズラ卵のサルモネラ汚染状況調査](https://mol.medicalonline.jp/library/journal/download?GoodsID=ee5jsofm/2016/003303/012&name=0160-0165j&UserID=133.50.190.185&base=jamas_pdf)'
$ws.Range("C26").Value = '未登録'

$ws.Range("A27").Value = '2011'
$ws.Range("B27").Value = '**農林水産省** <br> [各処理日の 1 番目・2 番目に処理されるブロイラー鶏群から製造された鶏肉のリステリア・モノサイトジェネス汚染状況の比較調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf)'
$ws.Range("C27").Value = '済'

$ws.Range("A28").Value = '2011'
$ws.Range("B28").Value = '**農林水産省** <br> [ブロイラー鶏群から製造された鶏肉のカンピロバクター汚染の季節変化調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf)'
$ws.Range("C28").Value = '済'

$ws.Range("A29").Value = '2011'
$ws.Range("B29").Value = '**農林水産省** <br> [ブロイラー鶏群から製造された鶏肉のサルモネラ汚染の季節変化調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf)'
$ws.Range("C29").Value = '済'

$ws.Range("A30").Value = '2011'
$ws.Range("B30").Value = '**農林水産省** <br> [豚の肝臓のカンピロバクター汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/butaniku/cam/02.html#24121)'
$ws.Range("C30").Value = '済'

$ws.Range("A31").Value = '2011'
$ws.Range("B31").Value = '**農林水産省** <br> [豚の肝臓のサルモネラ汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/butaniku/sal/02.html#24221)'
$ws.Range("C31").Value = '済'

$ws.Range("A32").Value = '2011'
$ws.Range("B32").Value = '**農林水産省** <br> [豚の肝臓のリステリア・モノサイトジェネス汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/butaniku/lis/02.html#24321)'
$ws.Range("C32").Value = '済'

$ws.Range("A33").Value = '2010-2013'
$ws.Range("B33").Value = '**東京都健康安全研究センター** <br> [東京都内に流通する牛内臓肉からの糞便系大腸菌群，ベロ毒素産生性大腸菌，_Campylobacter jejuni_/_coli_, _Salmonella_ および _Listeria monocytogenes_ 検出状況](https://www.jstage.jst.go.jp/article/jsfm/32/4/32_209/_pdf/-char/ja) <br>（日本食品微生物学会雑誌　2015）'
$ws.Range("C33").Value = '済'

$ws.Range("A34").Value = '2010'
$ws.Range("B34").Value = '**農林水産省** <br> [Prevalence and Characterization of Foodborne Pathogens in Dairy Cattle in the Eastern Part of Japan](https://www.jstage.jst.go.jp/article/jvms/75/4/75_12-0327/_pdf/-char/ja) <br> （Journal of Veterinary Medical Science 2013）'
$ws.Range("C34").Value = '済'

$ws.Range("A35").Value = '2010'
$ws.Range("B35").Value = '**農林水産省** <br> [出荷前後のブロイラー鶏群のカンピロバクター保有状況と、製造された鶏肉のカンピロバクター汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/keiniku_cam_14.html)'
$ws.Range("C35").Value = '済'

$ws.Range("A36").Value = '2010'
$ws.Range("B36").Value = '**農林水産省** <br> [ブロイラー鶏群から製造された中抜きと体及び鶏肉のカンピロバクター濃度調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf)'
$ws.Range("C36").Value = '済'

$ws.Range("A37").Value = '2010'
$ws.Range("B37").Value = '**農林水産省** <br> [ブロイラー鶏群から製造された中抜きと体及び鶏肉のサルモネラ濃度調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf)'
$ws.Range("C37").Value = '済'

$ws.Range("A38").Value = '2010'
$ws.Range("B38").Value = '**農林水産省** <br> [採卵鶏農場のサルモネラ保有状況・鶏卵のサルモネラ汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf) <br>  - [関連リンク](https://www.maff.go.jp/j/syouan/seisaku/kekka/keiran/keiran_sal_03.html#22113)'
$ws.Range("C38").Value = '済'

$ws.Range("A39").Value = '2010'
$ws.Range("B39").Value = '**福岡県保健環境研究所** <br> [平成22年度食品の食中毒菌汚染実態調査](https://www.fihes.pref.fukuoka.jp/nenpoh/np38/pdf/np38report1.pdf) <br> （福岡県保健環境研究所年報第38号, 66-67, 2011）'
$ws.Range("C39").Value = '未登録'

$ws.Range("A40").Value = '2009'
$ws.Range("B40").Value = '**農林水産省** <br> [ブロイラー鶏群から製造された鶏肉のカンピロバクター汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf) <br> - [関連リンク](https://www.maff.go.jp/j/syouan/seisaku/kekka/keiniku_cam_06.html)'
$ws.Range("C40").Value = '済'

$ws.Range("A41").Value = '2009'
$ws.Range("B41").Value = '**農林水産省** <br> [ブロイラー鶏群から製造された鶏肉のサルモネラ汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf)'
$ws.Range("C41").Value = '済'

$ws.Range("A42").Value = '2009'
$ws.Range("B42").Value = '**国立医薬品食品衛生研究所**<br>[_Campylobacter_ contamination in retail poultry meats and by-products in Japan: A literature survey](https://www.sciencedirect.com/science/article/pii/S0956713508002247)'
$ws.Range("C42").Value = '未登録'

$ws.Range("A43").Value = '2009'
$ws.Range("B43").Value = '**天使大学** <br> [Prevalence and Characteristics of _Listeria monocytogenes_ in Bovine Colostrum in Japan](https://www.sciencedirect.com/science/article/pii/S0362028X23052547)'
$ws.Range("C43").Value = '済'

$ws.Range("A44").Value = '2008-2018'
$ws.Range("B44").Value = '**厚生労働省** <br> [食品中の食中毒菌汚染実態調査](https://www.mhlw.go.jp/stf/seisakunitsuite/bunya/kenkou_iryou/shokuhin/syokuchu/01.html)'
$ws.Range("C44").Value = '済'

$ws.Range("A45").Value = '2007-2008'
$ws.Range("B45").Value = '**農林水産省** <br> [生食用野菜における腸管出血性大腸菌及びサルモネラの実態調査](https://www.maff.go.jp/j/syouan/nouan/kome/k_yasai/pdf/press.pdf)'
$ws.Range("C45").Value = '済'

$ws.Range("A46").Value = '2007'
$ws.Range("B46").Value = '**農林水産省** <br> [市販鶏卵のサルモネラ汚染状況調査](https://www.maff.go.jp/j/syouan/seisaku/kekka/attach/pdf/chikusan-1.pdf) <br> - [関連リンク](https://www.maff.go.jp/j/syouan/seisaku/kekka/keiran/keiran_sal_04.html#22121)'
$ws.Range("C46").Value = '済'

$ws.Range("A47").Value = '2007'
$ws.Range("B47").Value = '**帯広畜産大学** <br> [Isolation and characterization of _Listeria monocytogenes_ from commercial asazuke (Japanese light pickles)](https://www.sciencedirect.com/science/article/pii/S0168160510001686?via%3Dihub)'
$ws.Range("C47").Value = '未登録'

$ws.Range("A48").Value = '2007'
$ws.Range("B48").Value = '**日本獣医生命科学大学** <br> [Isolation of _Listeria monocytogenes_ from the Skin of Slaughtered Beef Cattle](https://agriknowledge.affrc.go.jp/RN/2030760055.pdf)'
$ws.Range("C48").Value = '未登録'

$ws.Range("A49").Value = '2006-2008'
$ws.Range("B49").Value = '**国立医薬品食品衛生研究所**<br>[The Occurrence of Listeria monocytogenes in Imported Ready-to-Eat Foods in Japan](https://www.jstage.jst.go.jp/article/jvms/74/3/74_11-0262/_pdf/-char/en) <br> (Journal of Veterinary Medical Science, 2012, Volume 74, Issue 3, Pages 373-375)'
$ws.Range("C49").Value = '済'

$ws.Range("A50").Value = '2006-2007'
$ws.Range("B50").Value = '**宮城県保健環境センター** <br> [芽物野菜等の食中毒菌汚染実態調査](https://www.pref.miyagi.jp/documents/1943/617283.pdf) <br> （宮城県保健環境センター年報　第26号, p.103-104,　2008）'
$ws.Range("C50").Value = '済'

$ws.Range("A51").Value = '2006'
$ws.Range("B51").Value = '**宮城県保健環境センター** <br> [_Listeria monocytogenes_ による ready-to-eat 食品の汚染実態](https://www.pref.miyagi.jp/documents/1979/210526.pdf) <br>（宮城県保健環境センター年報　第 25 号　2007）'
$ws.Range("C51").Value = '済'

$ws.Range("A52").Value = '2006'
$ws.Range("B52").Value = '**神戸大学** <br> [市販ミンチ肉における黄色ブドウ球菌汚染調査と分離株の性状](https://www.jstage.jst.go.jp/article/jsfm1994/23/4/23_4_217/_pdf/-char/ja) <br> （日本食品微生物学会雑誌, 23 (4), 217-222, 2006）'
$ws.Range("C52").Value = '未登録'

$ws.Range("A53").Value = '2004-2007'
$ws.Range("B53").Value = '**日本獣医生命科学大学** <br> [Sequence-Based Characterization of _Listeria monocytogenes_ Strains Isolated from Domestic Retail Meat in the Tokyo Metropolitan Area of Japan](https://www.jstage.jst.go.jp/article/yoken/71/5/71_JJID.2017.582/_pdf/-char/en)'
$ws.Range("C53").Value = '済'

$ws.Range("A54").Value = '2012'
$ws.Range("B54").Value = '**静岡県環境衛生科学研究所** <br> [Antibiotic Resistance in Bacterial Pathogens from Retail Raw Meats and Food-Producing Animals in Japan](https://www.sciencedirect.com/science/article/pii/S0362028X23039777?via%3Dihub) <br> (Journal of Food Protection, Volume 75, Issue 10, 1 October 2012, Pages 1774-1782)'
$ws.Range("C54").Value = '済'
